$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 now ends the anchor-word (A-H) table; the old "sc" anchor row (row 7)
# is dropped entirely, so clear its A:H cells (J:Q data for row 7 is kept/updated below).
$ws.Range("A7:H7").Clear()

# Updated values from the rerun (larger dataset)
$ws.Range('B3').Value = 0.8529411764705882
$ws.Range('C3').Value = 29
$ws.Range('D3').Value = 29
$ws.Range('H3').Value = 5
$ws.Range('J3').Value = 'happy'
$ws.Range('L3').Value = 26
$ws.Range('M3').Value = 26
$ws.Range('A4').Value = 'crisis'
$ws.Range('B4').Value = 0.6198630136986302
$ws.Range('C4').Value = 181
$ws.Range('D4').Value = 181
$ws.Range('H4').Value = 111
$ws.Range('J4').Value = 'interesting'
$ws.Range('K4').Value = 0.9696969696969697
$ws.Range('L4').Value = 32
$ws.Range('M4').Value = 32
$ws.Range('A5').Value = 'panic'
$ws.Range('B5').Value = 0.1705426356589147
$ws.Range('C5').Value = 88
$ws.Range('D5').Value = 88
$ws.Range('H5').Value = 428
$ws.Range('J5').Value = 'best'
$ws.Range('K5').Value = 0.9322033898305084
$ws.Range('L5').Value = 55
$ws.Range('M5').Value = 55
$ws.Range('Q5').Value = 4
$ws.Range('A6').Value = 'sc'
$ws.Range('B6').Value = 0.1428571428571428
$ws.Range('C6').Value = 27
$ws.Range('D6').Value = 27
$ws.Range('H6').Value = 162
$ws.Range('J6').Value = 'great'
$ws.Range('K6').Value = 0.875
$ws.Range('L6').Value = 98
$ws.Range('M6').Value = 98
$ws.Range('Q6').Value = 14
$ws.Range('J7').Value = 'love'
$ws.Range('K7').Value = 0.8695652173913043
$ws.Range('L7').Value = 40
$ws.Range('M7').Value = 40
$ws.Range('Q7').Value = 6
$ws.Range('J8').Value = 'special'
$ws.Range('K8').Value = 0.8333333333333334
$ws.Range('L8').Value = 30
$ws.Range('M8').Value = 30
$ws.Range('J9').Value = 'thanks'
$ws.Range('K9').Value = 0.8292682926829268
$ws.Range('L9').Value = 68
$ws.Range('M9').Value = 68
$ws.Range('Q9').Value = 14
$ws.Range('J10').Value = 'positive'
$ws.Range('K10').Value = 0.7931034482758621
$ws.Range('L10').Value = 46
$ws.Range('M10').Value = 46
$ws.Range('Q10').Value = 12
$ws.Range('J11').Value = 'thank'
$ws.Range('K11').Value = 0.78125
$ws.Range('L11').Value = 100
$ws.Range('M11').Value = 100
$ws.Range('Q11').Value = 28
$ws.Range('J12').Value = 'free'
$ws.Range('K12').Value = 0.7583333333333333
$ws.Range('L12').Value = 91
$ws.Range('M12').Value = 91
$ws.Range('Q12').Value = 29
$ws.Range('J13').Value = 'safe'
$ws.Range('K13').Value = 0.7394366197183099
$ws.Range('L13').Value = 105
$ws.Range('M13').Value = 105
$ws.Range('Q13').Value = 37
$ws.Range('J14').Value = 'safety'
$ws.Range('K14').Value = 0.7254901960784313
$ws.Range('L14').Value = 37
$ws.Range('M14').Value = 37
$ws.Range('Q14').Value = 14
$ws.Range('J15').Value = 'confidence'
$ws.Range('K15').Value = 0.7222222222222222
$ws.Range('L15').Value = 26
$ws.Range('M15').Value = 26
$ws.Range('Q15').Value = 10
$ws.Range('K16').Value = 0.7
$ws.Range('L16').Value = 112
$ws.Range('M16').Value = 112
$ws.Range('Q16').Value = 48
$ws.Range('J17').Value = 'support'
$ws.Range('K17').Value = 0.6792452830188679
$ws.Range('L17').Value = 72
$ws.Range('M17').Value = 72
$ws.Range('Q17').Value = 34
$ws.Range('J18').Value = 'better'
$ws.Range('K18').Value = 0.6190476190476191
$ws.Range('L18').Value = 39
$ws.Range('M18').Value = 39
$ws.Range('Q18').Value = 24
$ws.Range('J19').Value = 'relief'
$ws.Range('K19').Value = 0.6
$ws.Range('L19').Value = 30
$ws.Range('M19').Value = 30
$ws.Range('Q19').Value = 20
$ws.Range('J20').Value = 'well'
$ws.Range('K20').Value = 0.5851063829787234
$ws.Range('L20').Value = 55
$ws.Range('M20').Value = 55
$ws.Range('Q20').Value = 39
$ws.Range('J21').Value = 'fresh'
$ws.Range('K21').Value = 0.5833333333333334
$ws.Range('L21').Value = 28
$ws.Range('M21').Value = 28
$ws.Range('Q21').Value = 20
$ws.Range('J22').Value = 'hand'
$ws.Range('K22').Value = 0.5509138381201044
$ws.Range('L22').Value = 211
$ws.Range('M22').Value = 211
$ws.Range('Q22').Value = 172
$ws.Range('J23').Value = 'heroes'
$ws.Range('K23').Value = 0.5319148936170213
$ws.Range('L23').Value = 25
$ws.Range('M23').Value = 25
$ws.Range('Q23').Value = 22
$ws.Range('K24').Value = 0.4764705882352941
$ws.Range('L24').Value = 162
$ws.Range('M24').Value = 162
$ws.Range('Q24').Value = 178
$ws.Range('J25').Value = 'care'
$ws.Range('K25').Value = 0.449438202247191
$ws.Range('L25').Value = 40
$ws.Range('M25').Value = 40
$ws.Range('Q25').Value = 49
$ws.Range('J26').Value = 'help'
$ws.Range('K26').Value = 0.4305084745762712
$ws.Range('L26').Value = 127
$ws.Range('M26').Value = 127
$ws.Range('Q26').Value = 168
$ws.Range('J27').Value = 'protect'
$ws.Range('K27').Value = 0.3561643835616438
$ws.Range('L27').Value = 26
$ws.Range('M27').Value = 26
$ws.Range('Q27').Value = 47
$ws.Range('K28').Value = 0.3263598326359833
$ws.Range('L28').Value = 78
$ws.Range('M28').Value = 78
$ws.Range('Q28').Value = 161
